# Add a new "Category of data" column (C) to the data dictionary sheet,
# classifying each field as Numerical / Categorical / Binary / Doubt,
# plus a few notes in column D ("if slabs" / weekend note).
# Commit message: "Added Approx Category of data"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column C
$ws.Cells.Item(1, 3).Value = "To Del"

# Row_ID
$ws.Cells.Item(2, 3).Value = "Numerical"
# Days_taken_to_collect_data_after_publishing_article
$ws.Cells.Item(3, 3).Value = "Categorical"
# Title_Length
$ws.Cells.Item(4, 3).Value = "Numerical"
# Content_Length
$ws.Cells.Item(5, 3).Value = "Numerical"
# Rate_Unique_Words
$ws.Cells.Item(6, 3).Value = "Numerical"
# Rate_words_other_than_stop_words
$ws.Cells.Item(7, 3).Value = "Numerical"
# Rate_unique_words_other_than_stop_words
$ws.Cells.Item(8, 3).Value = "Categorical"
# Count_of_hyperlinks_in_content
$ws.Cells.Item(9, 3).Value = "Categorical"
# Count_of_miscellaneous_links
$ws.Cells.Item(10, 3).Value = "Categorical"
$ws.Cells.Item(10, 4).Value = "if slabs"
# Count_of_number_of_images
$ws.Cells.Item(11, 3).Value = "Categorical"
$ws.Cells.Item(11, 4).Value = "if slabs"
# Count_of_number_of_videos
$ws.Cells.Item(12, 3).Value = "Numerical"
# Average_word_length_in_content
$ws.Cells.Item(13, 3).Value = "Categorical"
# Count_Keywords_metadata
$ws.Cells.Item(14, 3).Value = "Binary"
# Lifestyle
$ws.Cells.Item(15, 3).Value = "Binary"
# Entertainment
$ws.Cells.Item(16, 3).Value = "Binary"
# Business
$ws.Cells.Item(17, 3).Value = "Binary"
# SocialMedia
$ws.Cells.Item(18, 3).Value = "Binary"
# Technology
$ws.Cells.Item(19, 3).Value = "Binary"
# World
$ws.Cells.Item(20, 3).Value = "Doubt"
# kw_min_min
$ws.Cells.Item(21, 3).Value = "Numerical"
# kw_max_min
$ws.Cells.Item(22, 3).Value = "Numerical"
# kw_avg_min
$ws.Cells.Item(23, 3).Value = "Doubt"
# kw_min_max
$ws.Cells.Item(24, 3).Value = "Numerical"
# kw_max_max
$ws.Cells.Item(25, 3).Value = "Numerical"
# kw_avg_max
$ws.Cells.Item(26, 3).Value = "Doubt"
# kw_min_avg
$ws.Cells.Item(27, 3).Value = "Numerical"
# kw_max_avg
$ws.Cells.Item(28, 3).Value = "Numerical"
# kw_avg_avg
$ws.Cells.Item(29, 3).Value = "Numerical"
# self_miscellaneous_min_shares
$ws.Cells.Item(30, 3).Value = "Numerical"
# self_miscellaneous_max_shares
$ws.Cells.Item(31, 3).Value = "Numerical"
# self_miscellaneous_avg_sharess
$ws.Cells.Item(32, 3).Value = "Binary"
# Monday
$ws.Cells.Item(33, 3).Value = "Binary"
# Tuesday
$ws.Cells.Item(34, 3).Value = "Binary"
# Wednesday
$ws.Cells.Item(35, 3).Value = "Binary"
# Thursday
$ws.Cells.Item(36, 3).Value = "Binary"
# Friday
$ws.Cells.Item(37, 3).Value = "Binary"
# Saturday
$ws.Cells.Item(38, 3).Value = "Binary"
# Sunday
$ws.Cells.Item(39, 3).Value = "Binary"
$ws.Cells.Item(39, 4).Value = "Used as per condition of Sat and Sun"
# Weekend
$ws.Cells.Item(40, 3).Value = "Numerical"
# miscellaneous 1
$ws.Cells.Item(41, 3).Value = "Numerical"
# miscellaneous 2
$ws.Cells.Item(42, 3).Value = "Numerical"
# miscellaneous 3
$ws.Cells.Item(43, 3).Value = "Numerical"
# miscellaneous 4
$ws.Cells.Item(44, 3).Value = "Numerical"
# miscellaneous 5
$ws.Cells.Item(45, 3).Value = "Numerical"
# global_subjectivity
$ws.Cells.Item(46, 3).Value = "Numerical"
# global_sentiment_polarity
$ws.Cells.Item(47, 3).Value = "Numerical"
# global_rate_positive_words
$ws.Cells.Item(48, 3).Value = "Numerical"
# global_rate_negative_words
$ws.Cells.Item(49, 3).Value = "Numerical"
# rate_positive_words
$ws.Cells.Item(50, 3).Value = "Numerical"
# rate_negative_words
$ws.Cells.Item(51, 3).Value = "Numerical"
# avg_positive_polarity
$ws.Cells.Item(52, 3).Value = "Categorical"
$ws.Cells.Item(52, 4).Value = "if slabs"
# min_positive_polarity
$ws.Cells.Item(53, 3).Value = "Categorical"
$ws.Cells.Item(53, 4).Value = "if slabs"
# max_positive_polarity
$ws.Cells.Item(54, 3).Value = "Categorical"
$ws.Cells.Item(54, 4).Value = "if slabs"
# avg_negative_polarity
$ws.Cells.Item(55, 3).Value = "Categorical"
$ws.Cells.Item(55, 4).Value = "if slabs"
# min_negative_polarity
$ws.Cells.Item(56, 3).Value = "Categorical"
$ws.Cells.Item(56, 4).Value = "if slabs"
# max_negative_polarity
$ws.Cells.Item(57, 3).Value = "Numerical"
# title_subjectivity
$ws.Cells.Item(58, 3).Value = "Numerical"
# title_sentiment_polarity
$ws.Cells.Item(59, 3).Value = "Numerical"
# abs_title_subjectivity
$ws.Cells.Item(60, 3).Value = "Numerical"
# abs_title_sentiment_polarity
$ws.Cells.Item(61, 3).Value = "Numerical"
# shares

# Highlight the new column with a yellow fill, like the author did
$ws.Columns.Item(3).Interior.Color = 65535

# Leave the whole new column selected, as it was when the workbook was saved
$null = $ws.Columns.Item(3).Select()
